# Reading data from excel and credentials from excel
# Adds a new "credentials" worksheet (after the existing "TextEditor" sheet)
# containing a username/password row, with the sample credential values
# styled the way a quick "looks like a hyperlink" / monospace callout would be.

$wb = $excel.ActiveWorkbook

# --- add the new sheet right after the current first/only sheet -----------
$firstSheet = $wb.Worksheets.Item(1)
$credSheet = $wb.Worksheets.Add([System.Type]::Missing, $firstSheet)
$credSheet.Name = "credentials"

# --- header + data ----------------------------------------------------------
$credSheet.Range("A1").Value = "username"
$credSheet.Range("B1").Value = "password"
$credSheet.Range("A2").Value = "ninjalinos@work.com"
$credSheet.Range("B2").Value = "sdet218920@"

# --- column widths, similar to the authored sheet ---------------------------
$credSheet.Columns.Item(1).ColumnWidth = 23.5703125
$credSheet.Columns.Item(2).ColumnWidth = 16.85546875

# --- formatting: username/password values get a small blue monospace font,
#     the username cell (email) additionally underlined like a hyperlink ----
$userFont = $credSheet.Range("A2").Font
$userFont.Name = "Consolas"
$userFont.Size = 10
$userFont.Color = 16711722
$userFont.Underline = $true

$passFont = $credSheet.Range("B2").Font
$passFont.Name = "Consolas"
$passFont.Size = 10
$passFont.Color = 16711722

# --- match page setup of the authored sheet ---------------------------------
$credSheet.PageSetup.Orientation = 1

# --- keep selection on B2, credentials becomes the active/visible tab ------
$credSheet.Range("B2").Select() | Out-Null
